$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - RandomForestRegressor values updated
$ws.Range("B3").Value = 8803897526765.172
$ws.Range("C3").Value = 5016646042772.495
$ws.Range("D3").Value = 4810874242826.667

# Row 4 - model renamed from GradientBoostingRegressor to DecisionTreeRegressor, values updated
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 1814373054278.681
$ws.Range("C4").Value = 1657753348927.594
$ws.Range("D4").Value = 1736063201603.138

# Row 5 - model renamed from AdaBoostRegressor to MLPRegressor, values updated
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 310830976286600.8
$ws.Range("C5").Value = 355315745924011.2
$ws.Range("D5").Value = 404108527231692.3
